$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "MalaTest"
$ws.Range("B11").Select()
